$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the whole "In the second section, Data Sources, ..." paragraph
#    (red bold call-out + trailing ".") that sits right before "Data sources:"
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "In the second section, Data Sources, list the database name and source, the independent variable and the predictors. Show ten rows of the database.^p",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 2) Remove the "In the third section..." / "List any predictor..." /
#    "Will you need to create dummy variables?" call-out paragraphs,
#    leaving two empty paragraphs behind.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "List any predictor or outcome variables that are categorical. ^p",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$found = $d.Content.Find.Execute(
    "In the third section, Model Details, use a scatter plot and explain how there is likely a linear relationship between the predictor variables and the outcome variable. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$found = $d.Content.Find.Execute(
    "Will you need to create dummy variables?",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 3) Append a new closing sentence to the final paragraph.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Even though the results for the two scatterplots weren" + [char]8217 + "t exactly what we were expecting, they still showed us some very interesting and useful relationships between the predictors and the output variable.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Even though the results for the two scatterplots weren" + [char]8217 + "t exactly what we were expecting, they still showed us some very interesting and useful relationships between the predictors and the output variable. Since we used numerical values we also didn" + [char]8217 + "t need to create any dummy variables.",
    2)

Write-Host "done"
